$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SalesMonthly")

# STORE NAME value (B2): "0 | All Stores" -> "GZQW2Y8 | "
$ws.Range("B2").Value = "GZQW2Y8 | "

# REPORTING PERIOD value (B4): "01-Aug-2023 To 31-Aug-2023" -> "01-Jul-2024 To 02-Feb-2025"
$ws.Range("B4").Value = "01-Jul-2024 To 02-Feb-2025"

# GGC TOTAL NET AMOUNT (B5): 32.07 -> 26.65
$ws.Range("B5").Value = 26.65

# Table header "month" label (A8): "Aug-2023" -> "Jan-2025"
$ws.Range("A8").Value = "Jan-2025"

# Data row 8 (table data row)
$ws.Range("B8").Value = 205
$ws.Range("E8").Value = 26.65
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 0
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 26.65
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 205
$ws.Range("L8").Value = 0

# Column B width adjustment
$ws.Range("B1").ColumnWidth = 27.607881546020508

# Force full precision recalculation setting
$excel.Application.CalculateFullRebuild()
